$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 225.53334
$ws.Range("I33").Value = 137.07692
$ws.Range("J33").Value = 800.5
$ws.Range("K33").Value = 137.07692
$ws.Range("L33").Value = 800.5
$ws.Range("M33").Value = 91.92308
$ws.Range("N33").Value = -1258.5
$ws.Range("H64").Value = 3914.7144
$ws.Range("I64").Value = 3900
$ws.Range("K64").Value = 3900
$ws.Range("M64").Value = -3652
$ws.Range("H67").Value = 3914.7144
$ws.Range("I67").Value = 3900
$ws.Range("K67").Value = 3900
$ws.Range("M67").Value = -3042
$ws.Range("H132").Value = 2235.5454
$ws.Range("I132").Value = 2159.1
$ws.Range("K132").Value = 6477.299999999999
$ws.Range("M132").Value = -3947.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3441.6099
$ws.Range("I32").Value = 3277.4285
$ws.Range("K32").Value = 3277.4285
$ws.Range("M32").Value = -2990.4285
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()
$ws.Range("H88").Value = 2183.923
$ws.Range("J88").Value = 2343.125
$ws.Range("L88").Value = 2343.125
$ws.Range("N88").Value = -3155.125
$ws.Range("H91").Value = 2183.923
$ws.Range("J91").Value = 2343.125
$ws.Range("L91").Value = 2343.125
$ws.Range("N91").Value = -5151.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9749.75
$ws.Range("I86").Value = 8000
$ws.Range("K86").Value = 8000
$ws.Range("M86").Value = -6877
$ws.Range("H89").Value = 9749.75
$ws.Range("I89").Value = 8000
$ws.Range("K89").Value = 40000
$ws.Range("M89").Value = -34384
$ws.Range("H99").Value = 1701.7
$ws.Range("I99").Value = 1536.1177
$ws.Range("J99").Value = 2640
$ws.Range("K99").Value = 1536.1177
$ws.Range("L99").Value = 2640
$ws.Range("M99").Value = -38.11770000000001
$ws.Range("N99").Value = -5636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1595.2632
$ws.Range("I7").Value = 927.53845
$ws.Range("K7").Value = 927.53845
$ws.Range("M7").Value = -814.53845
$ws.Range("H32").Value = 4956.3
$ws.Range("I32").Value = 2729.2222
$ws.Range("K32").Value = 2729.2222
$ws.Range("M32").Value = -2413.2222
$ws.Range("H74").Value = 25000
$ws.Range("J74").Value = 25000
$ws.Range("L74").Value = 25000
$ws.Range("N74").Value = -26748
$ws.Range("H77").Value = 25000
$ws.Range("J77").Value = 25000
$ws.Range("L77").Value = 75000
$ws.Range("N77").Value = -83736
$ws.Range("H97").Value = 16000
$ws.Range("J97").Value = 16000
$ws.Range("L97").Value = 16000
$ws.Range("N97").Value = -17982
$ws.Range("H134").Value = 1383.25
$ws.Range("I134").Value = 1399.1
$ws.Range("K134").Value = 4197.299999999999
$ws.Range("M134").Value = -1662.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 724.93335
$ws.Range("I34").Value = 574.1111
$ws.Range("K34").Value = 1722.3333
$ws.Range("M34").Value = -1638.3333
$ws.Range("H55").Value = 1090.909
$ws.Range("I55").Value = 777.7778
$ws.Range("J55").Value = 2500
$ws.Range("K55").Value = 2333.3334
$ws.Range("L55").Value = 7500
$ws.Range("M55").Value = -2156.3334
$ws.Range("N55").Value = -7854
$ws.Range("H75").Value = 6250
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 6250
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 18750
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -20746
$ws.Range("H78").Value = 6250
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 6250
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 56250
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -66234
$ws.Range("H95").Value = 13107.667
$ws.Range("I95").Value = 12999
$ws.Range("J95").Value = 13162
$ws.Range("K95").Value = 38997
$ws.Range("L95").Value = 39486
$ws.Range("M95").Value = -36938
$ws.Range("N95").Value = -43604
$ws.Range("H122").Value = 327.57144
$ws.Range("I122").Value = 152
$ws.Range("J122").Value = 397.8
$ws.Range("K122").Value = 1368
$ws.Range("L122").Value = 3580.2
$ws.Range("M122").Value = 1082
$ws.Range("N122").Value = -8480.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 917.3125
$ws.Range("I2").Value = 1698.1666
$ws.Range("K2").Value = 1698.1666
$ws.Range("M2").Value = -1585.1666
$ws.Range("H62").Value = 44000
$ws.Range("I62").Value = 44000
$ws.Range("K62").Value = 44000
$ws.Range("M62").Value = -43314
$ws.Range("H65").Value = 44000
$ws.Range("I65").Value = 44000
$ws.Range("K65").Value = 132000
$ws.Range("M65").Value = -128568
$ws.Range("H80").Value = 5245.909
$ws.Range("I80").Value = 2168.3333
$ws.Range("J80").Value = 6400
$ws.Range("K80").Value = 2168.3333
$ws.Range("L80").Value = 6400
$ws.Range("M80").Value = -1170.3333
$ws.Range("N80").Value = -8396
$ws.Range("H83").Value = 5245.909
$ws.Range("I83").Value = 2168.3333
$ws.Range("J83").Value = 6400
$ws.Range("K83").Value = 10841.6665
$ws.Range("L83").Value = 32000
$ws.Range("M83").Value = -5849.666499999999
$ws.Range("N83").Value = -41984
$ws.Range("H97").Value = 3621.5293
$ws.Range("I97").Value = 3557.9167
$ws.Range("J97").Value = 3774.2
$ws.Range("K97").Value = 3557.9167
$ws.Range("L97").Value = 3774.2
$ws.Range("M97").Value = -3061.9167
$ws.Range("N97").Value = -4766.2
$ws.Range("H102").Value = 1261.6154
$ws.Range("I102").Value = 1272.8182
$ws.Range("K102").Value = 1272.8182
$ws.Range("M102").Value = 349.1818000000001
$ws.Range("H126").Value = 8802.4
$ws.Range("I126").Value = 6337.3335
$ws.Range("J126").Value = 12500
$ws.Range("K126").Value = 19012.0005
$ws.Range("L126").Value = 37500
$ws.Range("M126").Value = -16542.0005
$ws.Range("N126").Value = -42440
$ws.Range("H138").Value = 110000
$ws.Range("J138").Value = 110000
$ws.Range("L138").Value = 110000
$ws.Range("N138").Value = -120280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("I40").Value = 2250
$ws.Range("K40").Value = 2250
$ws.Range("M40").Value = -2114
$ws.Range("H46").Value = 4044.5
$ws.Range("I46").Value = 3885.7856
$ws.Range("K46").Value = 3885.7856
$ws.Range("M46").Value = -3697.7856
$ws.Range("H55").Value = 438.42856
$ws.Range("I55").Value = 164.5
$ws.Range("J55").Value = 484.08334
$ws.Range("K55").Value = 164.5
$ws.Range("L55").Value = 484.08334
$ws.Range("M55").Value = 8.5
$ws.Range("N55").Value = -830.08334
$ws.Range("H61").Value = 2655.4285
$ws.Range("I61").Value = 2764.8333
$ws.Range("K61").Value = 2764.8333
$ws.Range("M61").Value = -2562.8333
$ws.Range("H63").Value = 49000
$ws.Range("I63").Value = 49000
$ws.Range("K63").Value = 49000
$ws.Range("M63").Value = -48251
$ws.Range("H66").Value = 49000
$ws.Range("I66").Value = 49000
$ws.Range("K66").Value = 147000
$ws.Range("M66").Value = -143256
$ws.Range("H68").Value = 6500
$ws.Range("J68").Value = 6500
$ws.Range("L68").Value = 6500
$ws.Range("N68").Value = -7998
$ws.Range("H71").Value = 6500
$ws.Range("J71").Value = 6500
$ws.Range("L71").Value = 32500
$ws.Range("N71").Value = -39988
$ws.Range("H93").Value = 1359.8
$ws.Range("J93").Value = 1449.75
$ws.Range("L93").Value = 1449.75
$ws.Range("N93").Value = -3945.75
$ws.Range("H99").Value = 15250
$ws.Range("I99").Value = 15250
$ws.Range("K99").Value = 15250
$ws.Range("M99").Value = -12255
$ws.Range("H113").Value = 2655.4285
$ws.Range("I113").Value = 2764.8333
$ws.Range("K113").Value = 2764.8333
$ws.Range("M113").Value = -594.8332999999998
$ws.Range("H122").Value = 6141.5
$ws.Range("I122").Value = 3283.3333
$ws.Range("J122").Value = 8999.666999999999
$ws.Range("K122").Value = 9849.999899999999
$ws.Range("L122").Value = 26999.001
$ws.Range("M122").Value = -7399.999899999999
$ws.Range("N122").Value = -31899.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1325.0869
$ws.Range("I100").Value = 1038.9231
$ws.Range("J100").Value = 1697.1
$ws.Range("K100").Value = 2077.8462
$ws.Range("L100").Value = 3394.2
$ws.Range("M100").Value = -1536.8462
$ws.Range("N100").Value = -4476.2
$ws.Range("H107").Value = 608.2
$ws.Range("I107").Value = 268.57144
$ws.Range("K107").Value = 805.71432
$ws.Range("M107").Value = 1114.28568
$ws.Range("H122").Value = 2689.0527
$ws.Range("I122").Value = 2784.5
$ws.Range("J122").Value = 2180
$ws.Range("K122").Value = 8353.5
$ws.Range("L122").Value = 6540
$ws.Range("M122").Value = -5903.5
$ws.Range("N122").Value = -11440
